$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (label unchanged) - update metric values
$ws.Range("B3").Value = 59478033414012.48
$ws.Range("C3").Value = 64931039046278.48
$ws.Range("D3").Value = 700147355913158.6

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update metric values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03626014045630512
$ws.Range("C4").Value = 0.03656582253708062
$ws.Range("D4").Value = 307778245403039.9

# Row 5: AdaBoostRegressor -> MLPRegressor, update metric values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 111807112420375
$ws.Range("C5").Value = 22819178645570.24
$ws.Range("D5").Value = 241537847368103.9
